$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.585.51"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.848.15"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -1.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.92"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3858"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.03"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07920"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9942"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.48"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "1.855.52"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.915"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.79"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06656"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001035"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "27.585.32"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.386"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.13"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.53"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.098"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.399"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.77"
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9738"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09394"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.573"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.283"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.346"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06017"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02223"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.306"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.181"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5901"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1862"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.29"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5586"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.19"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.897"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06688"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.70"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.054"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.11"
$ws.Range("E51").Value = "  -1.50%  "
